# Fa2020SemesterDates.xlsx — drop the "Spring semester" helper columns (G/H)
# that used to sit alongside the Fall dates, fill in the now-visible E31/E32
# weekday labels that the shared-formula fill had skipped, and leave the
# sheet scrolled/selected over the remaining A:B columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old Spring-semester date/weekday columns (G:H). Clearing
# contents (rather than deleting the columns outright) keeps the G column's
# formatting/style in place, matching the surviving empty <c s="1"/> cells.
$ws.Range("G1:H30").ClearContents() | Out-Null

# Rows 31/32 already had the weekday pattern (T/H) everywhere else in column
# E; fill in the two that were previously missing.
$ws.Range("E31").Value2 = "T"
$ws.Range("E32").Value2 = "H"

# Scroll down so row 23 is at the top of the view, and select A1:B48 (the
# data that remains relevant now that G:H are gone).
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1:B48").Select() | Out-Null
